$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the defined name (test_data_predictions_1 -> test_data_predictions) ---
foreach ($n in $wb.Names) {
    if ($n.Name -like "*test_data_predictions_1*") {
        $n.Name = "test_data_predictions"
    }
}

# --- Update header / label text (shared string pool content effectively swaps) ---
$ws.Range("E1").Value = "mme"
$ws.Range("D7").Value = "mmre"

# --- Row 2: predicted/error/relative-error values ---
$ws.Range("C2").Value = 318.92200000000003
$ws.Range("D2").Value = 88.221999999999994
$ws.Range("E2").Formula = "=ABS(C2-B2)/B2"

# --- Row 3: predicted/error/relative-error values (start of shared formula group) ---
$ws.Range("C3").Value = -255.19499999999999
$ws.Range("D3").Value = -327.19499999999999
$ws.Range("E3").Formula = "=ABS(C3-B3)/B3"

# --- Row 4 ---
$ws.Range("C4").Value = 168.95599999999999
$ws.Range("D4").Value = 38.655999999999999
$ws.Range("E4").Formula = "=ABS(C4-B4)/B4"

# --- Row 5 ---
$ws.Range("C5").Value = 326.00299999999999
$ws.Range("D5").Value = -10.297000000000001
$ws.Range("E5").Formula = "=ABS(C5-B5)/B5"

# --- Row 6 ---
$ws.Range("C6").Value = 190.108
$ws.Range("D6").Value = -96.891999999999996
$ws.Range("E6").Formula = "=ABS(C6-B6)/B6"

# --- Update the active selection shown in the sheet view ---
$ws.Range("J11").Select() | Out-Null

Write-Host "edit complete"
